$wb = $excel.ActiveWorkbook

# Sheet "Assets": F2 43101 -> 43104, F3 43103 -> 43101
$ws = $wb.Worksheets.Item("Assets")
$ws.Range("F2").Value = 43104
$ws.Range("F3").Value = 43101

# Sheet "Not assets": F2 43102 -> 43101, F3 43101 -> 43102
$ws = $wb.Worksheets.Item("Not assets")
$ws.Range("F2").Value = 43101
$ws.Range("F3").Value = 43102

# Sheet "Bank accounts": F2 43101 -> 43104, F3 43103 -> 43101
$ws = $wb.Worksheets.Item("Bank accounts")
$ws.Range("F2").Value = 43104
$ws.Range("F3").Value = 43101

# Sheet "Not bank accounts": F2 43102 -> 43101, F3 43101 -> 43102
$ws = $wb.Worksheets.Item("Not bank accounts")
$ws.Range("F2").Value = 43101
$ws.Range("F3").Value = 43102
